# Changes of 6th April 2022
# Refresh the FedEx PackageTrackNum (column C) / ShipmentTrackNum (column D)
# test values on Sheet1 with a new batch of tracking numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New PackageTrackNum values for rows 2..22 (column C).
$newTrackNums = @{
    2  = "320018252380"
    3  = "320018252391"
    4  = "320018252428"
    5  = "320018252461"
    6  = "320018252520"
    7  = "320018252564"
    8  = "320018252612"
    9  = "320018252656"
    10 = "320018252689"
    11 = "320018252704"
    12 = "320018252748"
    13 = "320018252760"
    14 = "320018252807"
    15 = "320018252829"
    16 = "320018252873"
    17 = "320018252910"
    18 = "320018252976"
    19 = "320018253012"
    20 = "320018253240"
    21 = "320018253284"
    22 = "320018253354"
}

# Rows where column D mirrors column C's PackageTrackNum value.
$mirrorRows = @(5, 6, 7, 13, 14, 15, 16, 17)

# The tracking numbers are all-digit strings, and a plain .Value assignment
# would otherwise be auto-converted to a number by Excel (same as typing
# them straight into the grid). Flip NumberFormat to Text on the cells we
# are about to overwrite before writing, then restore the normal style
# afterwards, exactly like it's done through the Excel UI.
$colCRange = $ws.Range("C2:C22")
$colCRange.NumberFormat = "@"
$ws.Range("D5:D7").NumberFormat = "@"
$ws.Range("D13:D17").NumberFormat = "@"

foreach ($row in $newTrackNums.Keys) {
    $value = $newTrackNums[$row]
    $ws.Cells.Item($row, 3).Value = $value
    if ($mirrorRows -contains $row) {
        $ws.Cells.Item($row, 4).Value = $value
    }
}

$colCRange.Style = "Normal"
$ws.Range("D5:D7").Style = "Normal"
$ws.Range("D13:D17").Style = "Normal"
